# Structure change - add usecase diagrams
# Rework the "Functional Requirements" table: remove the standalone
# "Payment Integration" requirement row, fold its related content into
# other rows, and reword a few descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional Requirements")

# --- Row 5 (FR-01, Inventory management) ---------------------------------
# Description is reworded to mention both Goods Received and Delivery Note.
$ws.Range("E5").Value = "Accountants shall be able to  manage inventory by creating a Goods Received/Delivery Note, when goods are imported into the warehouse or exported to resellers."
$ws.Rows.Item(5).RowHeight = 60

# --- Row 6 (FR-02, Order management) --------------------------------------
# Description stays the same text, just keep as-is (no textual change).
$ws.Range("E6").Value = "Accountants shall be able to process orders, track and update orders and payments status, generate invoices.  Resellers/customers should also be able to view the payment status of their orders."

# --- Row 7 (FR-03, Order Placement) ---------------------------------------
$ws.Range("C7").Value = "Order Placement"
$ws.Range("E7").Value = "Resellers/customers should be able to place an order for items by selecting the desired products, specifying the quantity, and choosing a payment method (Cash, bank transfer, Momo...)."

# --- Row 8 (FR-04, Order Status Tracking) ---------------------------------
$ws.Range("C8").Value = "Order Status Tracking"
$ws.Range("E8").Value = "Resellers/customers and accountants should be able to track the status of their orders, including knowing when the order has been processed, shipped, or delivered."
$ws.Rows.Item(8).RowHeight = 60

# --- Row 9 (FR-05, User Management) ---------------------------------------
$ws.Range("C9").Value = "User Management"
$ws.Range("E9").Value = "Admin shall be able to manage user profiles and accounts, assign user roles and permissions, control access to sensitive data, generating reports on users, …"

# --- Row 10 (FR-06, Product Management) -- unchanged ----------------------

# --- Row 11 (FR-07, Delivery Management) -----------------------------------
$ws.Range("C11").Value = "Delivery Management"
$ws.Range("E11").Value = "Accountants shall be able to create goods delivery note to deliver goods to resellers, update the status of orders as being transferred."

# --- Row 12 (FR-08, was Payment Integration) -- cleared out ---------------
$ws.Range("C12:E12").ClearContents()
$ws.Rows.Item(12).AutoFit()

# --- Row 13 (FR-09, Stock Reporting) ---------------------------------------
$ws.Range("C13").Value = "Stock Reporting"
$ws.Range("E13").Value = "The software should provide incoming/outgoing stock reports, allowing accountants to view stock levels, track inventory movements, and manage stock levels efficiently."

# --- Row 14 (FR-10, Sales Reporting) ---------------------------------------
$ws.Range("C14").Value = "Sales Reporting"
$ws.Range("E14").Value = "The software should provide sales reporting and revenue reports, allowing accountants to view best-selling products, track sales trends, and forecast sales performance."

# --- Update the active selection to match the new state -------------------
$ws.Activate()
$ws.Range("E9").Select()
